$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin price / 1h-volume refresh (values below are stored as plain text,
#     matching the feed's original inline-string cells) ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.938.86"
$ws.Range("E2").Value = "  +1.65%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.777.17"
$ws.Range("E3").Value = "  +1.76%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.19"
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4570"
$ws.Range("E7").Value = "  +2.76%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3596"
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07495"
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.95"
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.107"
$ws.Range("E11").Value = "  +1.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9997"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.88"
$ws.Range("E13").Value = "  +0.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.048"
$ws.Range("E14").Value = "  +0.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.231"
$ws.Range("E15").Value = "  +1.93%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.779.12"
$ws.Range("E16").Value = "  +1.43%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.71"
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001064"
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06432"
$ws.Range("E19").Value = "  +0.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9998"
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.23"
$ws.Range("E21").Value = "  +2.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.801"
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.960.69"
$ws.Range("E23").Value = "  +1.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.33"
$ws.Range("E24").Value = "  +1.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.086"
$ws.Range("E25").Value = "  -0.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.88"
$ws.Range("E26").Value = "  +0.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.29"
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.980.76"
$ws.Range("E28").Value = "  +1.46%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.224"
$ws.Range("E29").Value = "  +7.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.83"
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.127"
$ws.Range("E31").Value = "  +4.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09232"
$ws.Range("E32").Value = "  +2.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.668"
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.571"
$ws.Range("E34").Value = "  +0.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.90"
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02307"
$ws.Range("E36").Value = "  +1.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06172"
$ws.Range("E37").Value = "  +3.00%  "
$ws.Range("E38").Value = "  +0.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6336"
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.983"
$ws.Range("E40").Value = "  +1.04%  "
$ws.Range("E41").Value = "  -1.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.384"
$ws.Range("E42").Value = "  -0.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.906"
$ws.Range("E43").Value = "  +1.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.25"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "123.15"
$ws.Range("E47").Value = "  +1.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.963"
$ws.Range("E48").Value = "  +0.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06930"
$ws.Range("E49").Value = "  +1.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.138"
$ws.Range("E50").Value = "  -0.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.51"
$ws.Range("E51").Value = "  +0.41%  "

# --- Rows 45 and 46 swap places (Decentraland <-> PancakeSwap) with refreshed
#     price/volume figures ---
$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.742"
$ws.Range("E45").Value = "  +0.62%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5922"
$ws.Range("E46").Value = "  +1.10%  "
